$d = $word.ActiveDocument

# The last paragraph currently ends with " (NEAT algorithm)". We append a new
# list paragraph after it ("Point Cloud"), inheriting the same ListParagraph /
# numbered-list (numId 1, ilvl 0) formatting as the preceding bullet items.
$lastPara = $d.Paragraphs.Last
$tailRange = $lastPara.Range
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Collapse(0)
$newRange.InsertAfter("Point Cloud")
